# Update the cryptocurrency ranking snapshot (GitHub Actions refresh).
# Columns: A=idx B=Coin C=Link D=Price E=Volume(1h) F=Data G=Hora
# Only D (Price) and E (Volume 1h) values are refreshed for most rows;
# rows 7/8 additionally swap the Coin/Link between KuCoinToken and FTXToken.
#
# D/E columns hold text-formatted numeric strings (e.g. "310.79", "-1.05%")
# in the source workbook, so a leading apostrophe is used to force Excel to
# store them as text instead of auto-converting to numbers/percentages.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,4).Value = "'310.79"
$ws.Cells.Item(2,5).Value = "'-1.05%"

$ws.Cells.Item(3,4).Value = "'37.55"
$ws.Cells.Item(3,5).Value = "'-4.28%"

$ws.Cells.Item(4,4).Value = "'5.088"
$ws.Cells.Item(4,5).Value = "'-0.89%"

$ws.Cells.Item(5,4).Value = "'0.07769"
$ws.Cells.Item(5,5).Value = "'-4.52%"

$ws.Cells.Item(6,4).Value = "'4.346"
$ws.Cells.Item(6,5).Value = "'-3.41%"

$ws.Cells.Item(7,2).Value = "FTXToken"
$ws.Cells.Item(7,3).Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Cells.Item(7,4).Value = "'1.896"
$ws.Cells.Item(7,5).Value = "'-3.19%"

$ws.Cells.Item(8,2).Value = "KuCoinToken"
$ws.Cells.Item(8,3).Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Cells.Item(8,4).Value = "'8.206"
$ws.Cells.Item(8,5).Value = "'-1.02%"

$ws.Cells.Item(9,5).Value = "'-7.48%"

$ws.Cells.Item(10,4).Value = "'0.9173"
$ws.Cells.Item(10,5).Value = "'-2.26%"

$ws.Cells.Item(11,4).Value = "'0.1202"
$ws.Cells.Item(11,5).Value = "'-9.02%"

$ws.Cells.Item(12,4).Value = "'0.1925"
$ws.Cells.Item(12,5).Value = "'-2.22%"

$ws.Cells.Item(13,4).Value = "'0.08916"
$ws.Cells.Item(13,5).Value = "'-0.98%"

$ws.Cells.Item(14,4).Value = "'0.03414"
$ws.Cells.Item(14,5).Value = "'-2.20%"

$ws.Cells.Item(15,4).Value = "'0.09699"
$ws.Cells.Item(15,5).Value = "'-0.13%"

$ws.Cells.Item(16,4).Value = "'0.001368"
$ws.Cells.Item(16,5).Value = "'-2.77%"

$ws.Cells.Item(17,4).Value = "'0.005813"
$ws.Cells.Item(17,5).Value = "'-5.87%"

$ws.Cells.Item(18,4).Value = "'3.553"
$ws.Cells.Item(18,5).Value = "'-0.63%"

$ws.Cells.Item(19,4).Value = "'0.3393"
$ws.Cells.Item(19,5).Value = "'-2.10%"

$ws.Cells.Item(20,4).Value = "'0.1278"

$ws.Cells.Item(21,4).Value = "'5.032"
$ws.Cells.Item(21,5).Value = "'0.40%"

$ws.Cells.Item(23,5).Value = "'5,585.71%"

$ws.Cells.Item(24,4).Value = "'0.04370"
$ws.Cells.Item(24,5).Value = "'-0.21%"

$ws.Cells.Item(25,4).Value = "'0.001213"
$ws.Cells.Item(25,5).Value = "'-2.52%"

$ws.Cells.Item(26,4).Value = "'0.004255"
$ws.Cells.Item(26,5).Value = "'-10.03%"

$ws.Cells.Item(27,5).Value = "'-66.64%"

$ws.Cells.Item(39,4).Value = "'0.02113"
$ws.Cells.Item(39,5).Value = "'-4.23%"

$ws.Cells.Item(40,4).Value = "'0.04944"
$ws.Cells.Item(40,5).Value = "'-5.22%"

$ws.Cells.Item(41,4).Value = "'0.007639"
$ws.Cells.Item(41,5).Value = "'1.08%"

$ws.Cells.Item(42,5).Value = "'-4.43%"

$ws.Cells.Item(43,4).Value = "'0.1343"
$ws.Cells.Item(43,5).Value = "'-3.55%"

$ws.Cells.Item(44,4).Value = "'0.002058"
$ws.Cells.Item(44,5).Value = "'-2.13%"

$ws.Cells.Item(45,4).Value = "'0.009589"
$ws.Cells.Item(45,5).Value = "'5.15%"

$ws.Cells.Item(46,4).Value = "'0.00006684"
$ws.Cells.Item(46,5).Value = "'-2.00%"

$ws.Cells.Item(47,5).Value = "'-0.25%"

$ws.Cells.Item(48,4).Value = "'0.003040"
$ws.Cells.Item(48,5).Value = "'0.81%"

$ws.Cells.Item(50,5).Value = "'-0.25%"

$ws.Cells.Item(51,5).Value = "'-0.25%"
